# Adds the "anionic surfactants" ingredient family (12 new rows) to the
# HPC INCI worksheet, matching the upstream "Add files via upload" commit.
# New shared strings are created implicitly as each cell's .Value is set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new ingredient rows (35-46) for anionic surfactants section

# Row 35
$ws.Range('A35').Value = 'Sodium Lauryl Sulfate'
$ws.Range('B35').Value = 'lauril sulfato de sódio'
$ws.Range('C35').Value = 'laurilsulfato de sodio'
$ws.Range('D35').Value = 'Sodium Lauryl Sulfate'
$ws.Range('E35').Value = 'Sodium dodecyl sulfate'
$ws.Range('F35').Value = 'Sodium Lauryl Sulfate.png'
$ws.Range('G35').Value = 'anionic surfactants'
$ws.Range('H35').Value = 'widely used anionic surfactant known for its strong cleansing and foaming properties'

# Row 36
$ws.Range('A36').Value = 'Ammonium Lauryl Sulfate'
$ws.Range('B36').Value = 'lauril sulfato de amônio'
$ws.Range('C36').Value = 'lauril sulfato de amonio'
$ws.Range('D36').Value = 'Ammonium Lauryl Sulfate'
$ws.Range('E36').Value = 'Ammonium Dodecyl Sulfate'
$ws.Range('F36').Value = 'Ammonium Lauryl Sulfate.png'
$ws.Range('G36').Value = 'anionic surfactants'
$ws.Range('H36').Value = 'widely used anionic surfactant known for its strong cleansing and foaming properties'

# Row 37
$ws.Range('A37').Value = 'Triethanolamine Lauryl Sulfate'
$ws.Range('B37').Value = 'lauril sulfato de trietanolamina'
$ws.Range('C37').Value = 'lauril sulfato de trietanolamina'
$ws.Range('D37').Value = 'Triethanolamine Lauryl Sulfate'
$ws.Range('E37').Value = 'Tris(2-hydroxyethyl)ammonium dodecyl sulfate'
$ws.Range('F37').Value = 'Triethanolamine Lauryl Sulfate.png'
$ws.Range('G37').Value = 'anionic surfactants'
$ws.Range('H37').Value = 'anionic surfactant, acts as a cleansing and foaming agent, helping to remove dirt  while creating a rich, stable foam'

# Row 38
$ws.Range('A38').Value = 'Sodium Cetyl Sulfate'
$ws.Range('B38').Value = 'sulfato de cetila de sódio'
$ws.Range('C38').Value = 'sulfato de cetilo de sodio'
$ws.Range('D38').Value = 'Sodium Cetyl Sulfate'
$ws.Range('E38').Value = 'Sodium Hexadecyl Sulfate'
$ws.Range('F38').Value = 'Sodium Cetyl Sulfate.png'
$ws.Range('G38').Value = 'anionic surfactants'
$ws.Range('H38').Value = 'anionic surfactant , acts as a cleansing, emulsifying, and foaming agent, helping to remove dirt while stabilizing formulations'

# Row 39
$ws.Range('A39').Value = 'Polyoxyethylene Sodium Lauryl Ether Sulfate'
$ws.Range('B39').Value = 'lauril éter sulfato de sódio'
$ws.Range('C39').Value = 'lauril éter sulfato de socio'
$ws.Range('D39').Value = 'Sodium Laureth Sulfate'
$ws.Range('E39').Value = 'Sodium polyoxyethylene lauryl ether sulfate'
$ws.Range('F39').Value = 'Sodium lauryl ethoxysulfate.png'
$ws.Range('G39').Value = 'anionic surfactants'
$ws.Range('H39').Value = 'anionic surfactant, provides foaming, cleansing, and emulsifying properties, making it effective for removing dirt while maintaining a mild formulation'

# Row 40
$ws.Range('A40').Value = 'Sodium Lauroyl Methyl Taurate'
$ws.Range('B40').Value = 'lauril metil taurato de sódio'
$ws.Range('C40').Value = 'metil N-lauroil taurato de sodio'
$ws.Range('D40').Value = 'Sodium Lauroyl Methyl Taurate'
$ws.Range('E40').Value = 'Sodium 2-[methyl (1-oxododecyl)amino]ethanesulfonate'
$ws.Range('F40').Value = 'Sodium lauroyl methyl taurate.png'
$ws.Range('G40').Value = 'anionic surfactants'
$ws.Range('H40').Value = 'anionic surfactant,  provides mild cleansing, foaming, and emulsifying properties, making it effective for removing dirt while being gentle on the skin'

# Row 41
$ws.Range('A41').Value = 'Sodium Myristoyl Methyl Taurate'
$ws.Range('D41').Value = 'Sodium Myristoyl Methyl Taurate'
$ws.Range('E41').Value = 'Sodium 2-[methyl (1-oxotetradecyl)amino]ethanesulfonate'
$ws.Range('F41').Value = 'ethanesulfonic acid, 2-(methyl(1-oxotetradecyl)amino)-, sodium salt.png'
$ws.Range('G41').Value = 'anionic surfactants'
$ws.Range('H41').Value = 'anionic surfactant,  provides mild cleansing, foaming, and emulsifying properties, making it effective for removing dirt while being gentle on the skin'

# Row 42
$ws.Range('A42').Value = 'Sodium Lauroyl Methyl Alanine'
$ws.Range('D42').Value = 'Sodium Lauroyl Methyl Alanine'
$ws.Range('E42').Value = 'Sodium N-methyl-N-(1-oxododecyl)-β-alaninate'
$ws.Range('F42').Value = 'Sodium lauroyl methylaminopropionate.png'
$ws.Range('G42').Value = 'anionic surfactants'
$ws.Range('H42').Value = 'anionic surfactant,  provides mild cleansing, foaming, and emulsifying properties, making it effective for removing dirt while being gentle on the skin'

# Row 43
$ws.Range('A43').Value = 'Potassium Lauroyl Sarcosinate'
$ws.Range('B43').Value = 'potássio lauroil sarcosinato'
$ws.Range('C43').Value = 'potasio lauroil sarcosinato'
$ws.Range('D43').Value = 'Potassium Lauroyl Sarcosinate'
$ws.Range('E43').Value = 'Potassium salt of lauroyl sarcosine'
$ws.Range('F43').Value = 'Potassium Lauroyl Sarcosinate.png'
$ws.Range('G43').Value = 'anionic surfactants'
$ws.Range('H43').Value = 'anionic surfactant,  provides mild cleansing, foaming, and emulsifying properties, making it effective for removing dirt while being gentle on the skin'

# Row 44
$ws.Range('A44').Value = 'Triethanolamine lauroylsarcosinate'
$ws.Range('D44').Value = 'TEA Lauroyl Sarcosinate'
$ws.Range('E44').Value = 'N-Methyl-N-(1-oxododecyl)glycine, compound with 2,2'',2''''-nitrilotri(ethanol) (1:1)'
$ws.Range('F44').Value = 'Triethanolamine lauroylsarcosinate.png'
$ws.Range('H44').Value = 'anionic surfactant,  provides mild cleansing, foaming, and emulsifying properties, making it effective for removing dirt while being gentle on the skin'

# Row 45
$ws.Range('A45').Value = 'Sodium Alkyl (C14–16) Sulfonate'
$ws.Range('B45').Value = 'sódio C14-16 sulfonato de olefina'
$ws.Range('C45').Value = 'sulfonato de olefina (C14-16) sódico'
$ws.Range('D45').Value = 'Sodium C14-16 Olefin Sulfonate'
$ws.Range('E45').Value = 'Sulfonic acids, C14-16-alkane hydroxy and C14-16-alkene, sodium salts'
$ws.Range('F45').Value = 'Sodium c14 olefin sulfonate.png'
$ws.Range('G45').Value = 'anionic surfactants'
$ws.Range('H45').Value = 'anionic surfactant,  provides mild cleansing, foaming, and emulsifying properties, making it effective for removing dirt while being gentle on the skin'

# Row 46
$ws.Range('A46').Value = 'Dioctyl Sodium Sulfosuccinate'
$ws.Range('B46').Value = 'dioctilsulfossuccinato de sódio'
$ws.Range('C46').Value = 'sulfosuccinato de dioctilo, sal de sodio'
$ws.Range('D46').Value = 'Dioctyl Sodium Sulfosuccinate'
$ws.Range('E46').Value = 'Sodium 1,4-dioctoxy-1,4-dioxobutane-2-sulfonate'
$ws.Range('F46').Value = 'Di-n-octyl sodium sulfosuccinate.png'
$ws.Range('G46').Value = 'anionic surfactants'

# Restore the freeze-header view (top row frozen) and leave the same
# active cell the author ended up on after scrolling through the new rows.
$ws.Range('A2').Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.FreezePanes = $true
$ws.Range('D39').Select()
